{"js": "const body = context.document.body;\n\n// Original text -> replacement text pairs, derived from the diff.\nconst replacements = [\n  [\"2023-08-05 Saturday\", \"2023-08-06 Sunday\"],\n  [\"31\u00f72=15, 1\", \"92\u00f76=15, 2\"],\n  [\"58\u00f78=7, 2\", \"41\u00f73=13, 2\"],\n  [\"44\u00f77=6, 2\", \"53\u00f74=13, 1\"],\n  [\"43\u00f78=5, 3\", \"84\u00f77=12, 0\"],\n  [\"90\u00f77=12, 6\", \"38\u00f73=12, 2\"],\n  [\"20\u00f77=2, 6\", \"43\u00f77=6, 1\"],\n  [\"55\u00f74=13, 3\", \"94\u00f76=15, 4\"],\n  [\"44\u00f72=22, 0\", \"71\u00f74=17, 3\"],\n  [\"72\u00f76=12, 0\", \"84\u00f72=42, 0\"],\n  [\"87\u00f75=17, 2\", \"65\u00f78=8, 1\"],\n  [\"86\u00f77=12, 2\", \"73\u00f79=8, 1\"],\n  [\"63\u00f79=7, 0\", \"68\u00f74=17, 0\"],\n  [\"52\u00f76=8, 4\", \"67\u00f72=33, 1\"],\n  [\"52\u00f72=26, 0\", \"35\u00f77=5, 0\"],\n  [\"36\u00f74=9, 0\", \"80\u00f79=8, 8\"],\n  [\"81\u00f77=11, 4\", \"67\u00f78=8, 3\"],\n  [\"62\u00f72=31, 0\", \"88\u00f77=12, 4\"],\n  [\"97\u00f79=10, 7\", \"75\u00f79=8, 3\"],\n  [\"11\u00f76=1, 5\", \"61\u00f77=8, 5\"],\n  [\"20\u00f76=3, 2\", \"96\u00f75=19, 1\"],\n  [\"20\u00f72=10, 0\", \"33\u00f74=8, 1\"],\n  [\"39\u00f79=4, 3\", \"19\u00f79=2, 1\"],\n  [\"90\u00f76=15, 0\", \"74\u00f77=10, 4\"],\n  [\"49\u00f73=16, 1\", \"72\u00f76=12, 0\"],\n  [\"88\u00f73=29, 1\", \"53\u00f72=26, 1\"],\n];\n\n// Phase 1: search for every original string *before* any edits are made.\n// This avoids a just-inserted replacement text being accidentally matched\n// by a later search (one of the new strings also happens to be an old string\n// found elsewhere in the document).\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"text\"));\nawait context.sync();\n\n// Phase 2: replace each matched range with its new text.\nsearchResults.forEach((result, i) => {\n  const [oldText, newText] = replacements[i];\n  if (result.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\", found ${result.items.length}`\n    );\n  }\n  result.items[0].insertText(newText, Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# Applies the diff: updates the date heading and the 25 division-problem\n# table cells. Each (row, column) target is addressed directly via the\n# Tables/Cell object model, so there is no ambiguity even though one of the\n# new values (\"72\u00f76=12, 0\") also happens to be an old value used elsewhere\n# in the table.\n$d = $word.ActiveDocument\n\n# Paragraph/cell Range.Text includes trailing mark characters (paragraph\n# mark \\r, cell mark \\a, ...); strip those before comparing against the\n# plain text we expect to find.\nfunction Get-PlainText($rangeText) {\n    return $rangeText.TrimEnd([char]13, [char]7)\n}\n\n# --- Date heading (first paragraph, above the table) ---\n$dateParagraph = $d.Paragraphs.Item(1)\n$expectedDate = '2023-08-05 Saturday'\n$actualDate = Get-PlainText $dateParagraph.Range.Text\nif ($actualDate -ne $expectedDate) {\n    throw \"Unexpected date heading text: [$actualDate]\"\n}\n$dateParagraph.Range.Text = '2023-08-06 Sunday'\n\n# --- Division-problem table cells ---\n$table = $d.Tables.Item(1)\n\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; OldText = '31\u00f72=15, 1'; NewText = '92\u00f76=15, 2' }\n    @{ Row = 1; Col = 2; OldText = '58\u00f78=7, 2'; NewText = '41\u00f73=13, 2' }\n    @{ Row = 1; Col = 3; OldText = '44\u00f77=6, 2'; NewText = '53\u00f74=13, 1' }\n    @{ Row = 1; Col = 4; OldText = '43\u00f78=5, 3'; NewText = '84\u00f77=12, 0' }\n    @{ Row = 1; Col = 5; OldText = '90\u00f77=12, 6'; NewText = '38\u00f73=12, 2' }\n    @{ Row = 5; Col = 1; OldText = '20\u00f77=2, 6'; NewText = '43\u00f77=6, 1' }\n    @{ Row = 5; Col = 2; OldText = '55\u00f74=13, 3'; NewText = '94\u00f76=15, 4' }\n    @{ Row = 5; Col = 3; OldText = '44\u00f72=22, 0'; NewText = '71\u00f74=17, 3' }\n    @{ Row = 5; Col = 4; OldText = '72\u00f76=12, 0'; NewText = '84\u00f72=42, 0' }\n    @{ Row = 5; Col = 5; OldText = '87\u00f75=17, 2'; NewText = '65\u00f78=8, 1' }\n    @{ Row = 9; Col = 1; OldText = '86\u00f77=12, 2'; NewText = '73\u00f79=8, 1' }\n    @{ Row = 9; Col = 2; OldText = '63\u00f79=7, 0'; NewText = '68\u00f74=17, 0' }\n    @{ Row = 9; Col = 3; OldText = '52\u00f76=8, 4'; NewText = '67\u00f72=33, 1' }\n    @{ Row = 9; Col = 4; OldText = '52\u00f72=26, 0'; NewText = '35\u00f77=5, 0' }\n    @{ Row = 9; Col = 5; OldText = '36\u00f74=9, 0'; NewText = '80\u00f79=8, 8' }\n    @{ Row = 13; Col = 1; OldText = '81\u00f77=11, 4'; NewText = '67\u00f78=8, 3' }\n    @{ Row = 13; Col = 2; OldText = '62\u00f72=31, 0'; NewText = '88\u00f77=12, 4' }\n    @{ Row = 13; Col = 3; OldText = '97\u00f79=10, 7'; NewText = '75\u00f79=8, 3' }\n    @{ Row = 13; Col = 4; OldText = '11\u00f76=1, 5'; NewText = '61\u00f77=8, 5' }\n    @{ Row = 13; Col = 5; OldText = '20\u00f76=3, 2'; NewText = '96\u00f75=19, 1' }\n    @{ Row = 17; Col = 1; OldText = '20\u00f72=10, 0'; NewText = '33\u00f74=8, 1' }\n    @{ Row = 17; Col = 2; OldText = '39\u00f79=4, 3'; NewText = '19\u00f79=2, 1' }\n    @{ Row = 17; Col = 3; OldText = '90\u00f76=15, 0'; NewText = '74\u00f77=10, 4' }\n    @{ Row = 17; Col = 4; OldText = '49\u00f73=16, 1'; NewText = '72\u00f76=12, 0' }\n    @{ Row = 17; Col = 5; OldText = '88\u00f73=29, 1'; NewText = '53\u00f72=26, 1' }\n)\n\nforeach ($update in $cellUpdates) {\n    $cell = $table.Cell($update.Row, $update.Col)\n    $actual = Get-PlainText $cell.Range.Text\n    if ($actual -ne $update.OldText) {\n        throw \"Unexpected text in cell ($($update.Row),$($update.Col)): [$actual]\"\n    }\n    $cell.Range.Text = $update.NewText\n}\n\n"}
